# Partner inventory invoice annexure - add Booking ID / Rate columns
# (CRM-1010: change invoice generation process for returned defective inventory)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column C, shifting C..G to D..H.
#    Default Insert() copies formatting from the LEFT column; we want the
#    new column to inherit the formatting that the (old) column C had, so
#    immediately after the insert we paste formats from the new column D
#    (which now holds what used to be column C) back onto the new column C.
# ---------------------------------------------------------------------------
$ws.Columns("C").Insert()
$ws.Range("D1:D18").Copy()
$ws.Range("C1:C18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Restore the explicit column width for the (now duplicated) column C - the
# insert only carried the width onto column D.
$ws.Columns("C").ColumnWidth = $ws.Columns("D").ColumnWidth

# ---------------------------------------------------------------------------
# 2. The header/table row (16-17) was hand edited rather than just shifted:
#    the B:C and E:F merges stay put, "Quantity" moves from column D to G
#    (as literal new content - "Booking ID" takes over D) and a new "Rate"
#    column appears in H. Undo the automatic merge growth from step 1 and
#    rebuild the row explicitly.
# ---------------------------------------------------------------------------
$ws.Range("B16:D16").UnMerge()
$ws.Range("B17:D17").UnMerge()
$ws.Range("F16:H16").UnMerge()
$ws.Range("F17:H17").UnMerge()

# Clear the stray cells left over from the shift so only the intended cells
# carry content.
$ws.Range("C16:H17").ClearContents()

# -- Row 16 (headers) --
$ws.Range("B16").Value = "Reference Invoice"
$ws.Range("D16").Value = "Booking ID"
$ws.Range("E16").Value = "Part Number"
$ws.Range("G16").Value = "Quantity"
$ws.Range("H16").Value = "Rate"

# -- Row 17 (placeholder tokens) --
$ws.Range("B17").Value = "{booking:incoming_invoice_id}"
$ws.Range("D17").Value = "{booking:booking_id}"
$ws.Range("E17").Value = "{booking:part_number}"
$ws.Range("G17").Value = "{booking:qty}"
$ws.Range("H17").Value = "{booking:rate}"

# Re-apply consistent formatting across the header row (grey fill, bold
# Cambria, left/bottom aligned) by copying from B16's original formatting.
$ws.Range("B16").Copy()
$ws.Range("D16,E16,G16,H16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Row 17 formatting: B/E/G share one look (regular Cambria, bottom aligned),
# D/H share another (regular Cambria, vertically centered).
$ws.Range("B17").Copy()
$ws.Range("E17,G17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("D17").VerticalAlignment = -4108
$ws.Range("H17").VerticalAlignment = -4108
$ws.Range("D17").HorizontalAlignment = -4131
$ws.Range("H17").HorizontalAlignment = -4131

# Re-create the merges the way the target layout wants them.
$ws.Range("B16:C16").Merge()
$ws.Range("B17:C17").Merge()
$ws.Range("E16:F16").Merge()
$ws.Range("E17:F17").Merge()

Write-Host "done"
